$wb = $excel.ActiveWorkbook

# "secondary" is sheet1, "primary" is sheet2
$ws1 = $wb.Worksheets.Item("secondary")

# Fix the bad test fixture values on the "secondary" sheet
$ws1.Range("C4").Value = "STTC-03"
$ws1.Range("C5").Value = "STTC-04"

# Make "secondary" the active sheet (it was incorrectly "primary" before)
# and select C6 there, matching the corrected fixture state.
$ws1.Activate()
$ws1.Range("C6").Select()
